$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Convert column A (rows 2-22) from numeric IDs to text IDs ---
# Excel auto-coerces a numeric-looking string assigned to .Value back into a
# number, so force the cells to Text format first (this is what makes the
# written value persist as a string instead of a number), write the values,
# then restore the original ("Normal") style so no stray number-format is
# left behind on the cells.
$idValues = @{
    2  = "101"
    3  = "123"
    4  = "4"
    5  = "5"
    6  = "6"
    7  = "7"
    8  = "8"
    9  = "9"
    10 = "10"
    11 = "11"
    12 = "12"
    13 = "13"
    14 = "14"
    15 = "15"
    16 = "16"
    17 = "17"
    18 = "18"
    19 = "19"
    20 = "20"
    21 = "21"
    22 = "22"
}

$ws.Range("A2:A22").NumberFormat = "@"
foreach ($r in $idValues.Keys) {
    $ws.Cells.Item($r, 1).Value = $idValues[$r]
}
$ws.Range("A2:A22").Style = "Normal"

# --- Append new client rows (23-26) ---
$newRows = @(
    @{ A = "59595959595"; B = "Nathalia Cunha"; C = 234; D = "Rua Rio Tocantins 845"; E = 0 },
    @{ A = "11568973738"; B = "Jurandir Silva"; C = 34;  D = "Rua B";                 E = 0 },
    @{ A = "11111111111"; B = "aaaaaaaa";       C = 34;  D = "aaaaaaaa";              E = 0 },
    @{ A = "1647894877";  B = "Gustavo Borges";  C = 34;  D = "Rua A";                 E = 0 }
)

$ws.Range("A23:A26").NumberFormat = "@"
$row = 23
foreach ($entry in $newRows) {
    $ws.Cells.Item($row, 1).Value = $entry.A
    $ws.Cells.Item($row, 2).Value = $entry.B
    $ws.Cells.Item($row, 3).Value = $entry.C
    $ws.Cells.Item($row, 4).Value = $entry.D
    $ws.Cells.Item($row, 5).Value = $entry.E
    $row++
}
$ws.Range("A23:A26").Style = "Normal"
